$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6 (replacing the old 2000-2004 data)
$years = @("2010年", "2011年", "2012年", "2013年", "2014年")
$values = @(308622, 168745, 157807, 113439, 14377)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Delete old rows 7-16 which are no longer present
$ws.Range("A7:B16").EntireRow.Delete()
